$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2321.6785
$ws.Range("I137").Value = 2337.15
$ws.Range("J137").Value = 2283
$ws.Range("K137").Value = 7011.450000000001
$ws.Range("L137").Value = 6849
$ws.Range("M137").Value = -4461.450000000001
$ws.Range("N137").Value = -11949
$ws.Range("H141").Value = 6293.8335
$ws.Range("I141").Value = 5465.1
$ws.Range("K141").Value = 16395.3
$ws.Range("M141").Value = -11215.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3477.2222
$ws.Range("I32").Value = 2097.1667
$ws.Range("K32").Value = 2097.1667
$ws.Range("M32").Value = -1810.1667
$ws.Range("H61").Value = 1677.3334
$ws.Range("I61").Value = 1012.9
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 1012.9
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -800.9
$ws.Range("N61").Value = -5423.5
$ws.Range("H96").Value = 19373.75
$ws.Range("J96").Value = 19373.75
$ws.Range("L96").Value = 19373.75
$ws.Range("N96").Value = -24865.75
$ws.Range("H122").Value = 7301.1304
$ws.Range("I122").Value = 7408.75
$ws.Range("K122").Value = 22226.25
$ws.Range("M122").Value = -19776.25
$ws.Range("H136").Value = 1677.3334
$ws.Range("I136").Value = 1012.9
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 3038.7
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -488.6999999999998
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1311.2174
$ws.Range("I20").Value = 1419.3846
$ws.Range("J20").Value = 1170.6
$ws.Range("K20").Value = 1419.3846
$ws.Range("L20").Value = 1170.6
$ws.Range("M20").Value = -1172.3846
$ws.Range("N20").Value = -1664.6
$ws.Range("H99").Value = 76924820
$ws.Range("I99").Value = 90910880
$ws.Range("K99").Value = 90910880
$ws.Range("M99").Value = -90909382
$ws.Range("H107").Value = 10271.637
$ws.Range("I107").Value = 9682.263000000001
$ws.Range("J107").Value = 14004.333
$ws.Range("K107").Value = 9682.263000000001
$ws.Range("L107").Value = 14004.333
$ws.Range("M107").Value = -7762.263000000001
$ws.Range("N107").Value = -17844.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1666.6923
$ws.Range("I31").Value = 1388.9166
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1388.9166
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1093.9166
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 1666.6923
$ws.Range("I34").Value = 1388.9166
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1388.9166
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1186.9166
$ws.Range("N34").Value = -5404
$ws.Range("H58").Value = 6027.143
$ws.Range("J58").Value = 2497
$ws.Range("L58").Value = 2497
$ws.Range("N58").Value = -2903
$ws.Range("H99").Value = 4175.1665
$ws.Range("I99").Value = 2327.75
$ws.Range("K99").Value = 2327.75
$ws.Range("M99").Value = -829.75
$ws.Range("H126").Value = 4175.1665
$ws.Range("I126").Value = 2327.75
$ws.Range("K126").Value = 6983.25
$ws.Range("M126").Value = -4513.25
$ws.Range("H132").Value = 2411.2593
$ws.Range("I132").Value = 1385.25
$ws.Range("K132").Value = 4155.75
$ws.Range("M132").Value = -1625.75
$ws.Range("H136").Value = 6027.143
$ws.Range("J136").Value = 2497
$ws.Range("L136").Value = 7491
$ws.Range("N136").Value = -12591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 21458452
$ws.Range("J2").Value = 47619176
$ws.Range("L2").Value = 285715056
$ws.Range("N2").Value = -285715282
$ws.Range("H37").Value = 80402.71000000001
$ws.Range("J37").Value = 80402.71000000001
$ws.Range("L37").Value = 241208.13
$ws.Range("N37").Value = -241432.13
$ws.Range("H46").Value = 96682400
$ws.Range("I46").Value = 55555610
$ws.Range("J46").Value = 105821690
$ws.Range("K46").Value = 166666830
$ws.Range("L46").Value = 317465070
$ws.Range("M46").Value = -166666739
$ws.Range("N46").Value = -317465252
$ws.Range("H68").Value = 999.25
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 999.25
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H86").Value = 507
$ws.Range("I86").Value = 508.2857
$ws.Range("K86").Value = 1524.8571
$ws.Range("M86").Value = -338.8571000000002
$ws.Range("H89").Value = 507
$ws.Range("I89").Value = 508.2857
$ws.Range("K89").Value = 4574.571300000001
$ws.Range("M89").Value = 1353.428699999999
$ws.Range("H107").Value = 1149.5
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 1199
$ws.Range("K107").Value = 3300
$ws.Range("L107").Value = 3597
$ws.Range("M107").Value = -1380
$ws.Range("N107").Value = -7437
$ws.Range("H122").Value = 980.4400000000001
$ws.Range("I122").Value = 2276.6
$ws.Range("J122").Value = 656.4
$ws.Range("K122").Value = 20489.4
$ws.Range("L122").Value = 5907.599999999999
$ws.Range("M122").Value = -18039.4
$ws.Range("N122").Value = -10807.6
$ws.Range("H136").Value = 2380
$ws.Range("I136").Value = 2380
$ws.Range("K136").Value = 7140
$ws.Range("M136").Value = -2040
$ws.Range("H137").Value = 10000
$ws.Range("I137").Value = 10000
$ws.Range("K137").Value = 30000
$ws.Range("M137").Value = -24900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8874.5
$ws.Range("J92").Value = 8499.333000000001
$ws.Range("L92").Value = 8499.333000000001
$ws.Range("N92").Value = -12243.333
$ws.Range("H113").Value = 1522.5454
$ws.Range("I113").Value = 1484
$ws.Range("K113").Value = 1484
$ws.Range("M113").Value = 686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6255.2915
$ws.Range("I40").Value = 6029.3076
$ws.Range("J40").Value = 6522.364
$ws.Range("K40").Value = 6029.3076
$ws.Range("L40").Value = 6522.364
$ws.Range("M40").Value = -5893.3076
$ws.Range("N40").Value = -6794.364
$ws.Range("H46").Value = 1249.5
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812
$ws.Range("H55").Value = 589.8461
$ws.Range("J55").Value = 1030
$ws.Range("L55").Value = 1030
$ws.Range("N55").Value = -1376
$ws.Range("H93").Value = 4459.9
$ws.Range("I93").Value = 2491.7856
$ws.Range("J93").Value = 9052.166999999999
$ws.Range("K93").Value = 2491.7856
$ws.Range("L93").Value = 9052.166999999999
$ws.Range("M93").Value = -1243.7856
$ws.Range("N93").Value = -11548.167
$ws.Range("H122").Value = 5161.558
$ws.Range("I122").Value = 4131.125
$ws.Range("J122").Value = 8159.1816
$ws.Range("K122").Value = 12393.375
$ws.Range("L122").Value = 24477.5448
$ws.Range("M122").Value = -9943.375
$ws.Range("N122").Value = -29377.5448
$ws.Range("H123").Value = 33332.555
$ws.Range("J123").Value = 33332.555
$ws.Range("L123").Value = 33332.555
$ws.Range("N123").Value = -43132.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 24074
$ws.Range("J86").Value = 24074
$ws.Range("L86").Value = 24074
$ws.Range("N86").Value = -26320
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 24074
$ws.Range("J89").Value = 24074
$ws.Range("L89").Value = 120370
$ws.Range("N89").Value = -131602
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H126").Value = 2256.2856
$ws.Range("I126").Value = 1965
$ws.Range("K126").Value = 5895
$ws.Range("M126").Value = -3425
$ws.Range("H132").Value = 1023.0909
$ws.Range("I132").Value = 922.375
$ws.Range("K132").Value = 2767.125
$ws.Range("M132").Value = -237.125
